$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "25÷8=" "68÷3="
Replace-Text "85÷8=" "92÷8="
Replace-Text "74÷9=" "39÷2="
Replace-Text "40÷6=" "41÷7="
Replace-Text "89÷5=" "22÷7="
Replace-Text "22÷6=" "45÷5="
Replace-Text "45÷9=" "46÷5="
Replace-Text "23÷2=" "20÷5="
Replace-Text "70÷9=" "53÷5="
Replace-Text "34÷2=" "83÷3="
Replace-Text "22÷9=" "14÷5="
Replace-Text "11÷6=" "54÷3="
Replace-Text "68÷8=" "90÷4="
Replace-Text "21÷8=" "97÷7="
Replace-Text "75÷6=" "93÷8="
Replace-Text "17÷6=" "84÷9="
Replace-Text "18÷3=" "17÷5="
Replace-Text "48÷3=" "92÷4="
Replace-Text "59÷7=" "73÷8="
Replace-Text "19÷2=" "11÷8="
Replace-Text "36÷5=" "14÷5="
Replace-Text "81÷6=" "71÷5="
Replace-Text "76÷2=" "87÷3="
Replace-Text "60÷5=" "61÷9="
Replace-Text "96÷3=" "65÷6="
